$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.816.06'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.870.61'
$ws.Range('E3').Value = '  -2.17%  '
$ws.Range('E4').Value = '  -0.92%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.95'
$ws.Range('E5').Value = '  -3.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.675'
$ws.Range('E6').Value = '  -6.13%  '
$ws.Range('E7').Value = '  -0.97%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.13'
$ws.Range('E8').Value = '  +3.45%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.341'
$ws.Range('E9').Value = '  -5.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0733'
$ws.Range('E10').Value = '  -2.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0966'
$ws.Range('E11').Value = '  -2.56%  '
$ws.Range('B12').Value = 'Chainlink'
$ws.Range('C12').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '12.83'
$ws.Range('E12').Value = '  +1.70%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.141.49'
$ws.Range('E13').Value = '  -2.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.710'
$ws.Range('E14').Value = '  -1.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.866.38'
$ws.Range('E15').Value = '  -2.34%  '
$ws.Range('E16').Value = '  -1.96%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '34.806.17'
$ws.Range('E17').Value = '  -1.47%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '72.20'
$ws.Range('E18').Value = '  -2.67%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0809'
$ws.Range('E19').Value = '  -4.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '242.51'
$ws.Range('E20').Value = '  -0.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.55'
$ws.Range('E21').Value = '  -3.30%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.88'
$ws.Range('E22').Value = '  -3.94%  '
$ws.Range('E23').Value = '  -0.93%  '
$ws.Range('E24').Value = '  +4.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.14'
$ws.Range('E25').Value = '  -12.41%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.34'
$ws.Range('E26').Value = '  -2.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.33'
$ws.Range('E27').Value = '  -3.37%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.02'
$ws.Range('E28').Value = '  -3.77%  '
$ws.Range('E29').Value = '  -5.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.128.50'
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('E31').Value = '  +6.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.17'
$ws.Range('E32').Value = '  -4.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0569'
$ws.Range('E33').Value = '  -2.09%  '
$ws.Range('E34').Value = '  -0.99%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.12'
$ws.Range('E35').Value = '  -1.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.828'
$ws.Range('E36').Value = '  -9.46%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.94'
$ws.Range('E37').Value = '  -4.47%  '
$ws.Range('E38').Value = '  -25.46%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '97.43'
$ws.Range('E39').Value = '  +0.43%  '
$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '16.87'
$ws.Range('E40').Value = '  -2.69%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0664'
$ws.Range('E41').Value = '  +2.75%  '
$ws.Range('E42').Value = '  -3.80%  '
$ws.Range('E43').Value = '  -3.66%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.281.25'
$ws.Range('E44').Value = '  -4.25%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0821'
$ws.Range('E45').Value = '  +11.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.30'
$ws.Range('E46').Value = '  -5.04%  '
$ws.Range('E47').Value = '  -1.22%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.73'
$ws.Range('E48').Value = '  -1.50%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '11.79'
$ws.Range('E49').Value = '  -4.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.28'
$ws.Range('E50').Value = '  -7.19%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '42.57'
$ws.Range('E51').Value = '  -5.70%  '
